# Exámenes Ordinario/Extraordinario 30 mayo de 2024 - Lap HP
# Adds two new exam-date columns (6-May-2024 and 13-May-2024) to the
# "Asistencias" sheet, pushing the "Faltas" totals column from L to N,
# and switches the active sheet/selection to "Asistencias".

$wb = $excel.ActiveWorkbook

$concentrado = $wb.Worksheets.Item("Concentrado")
$asistencias = $wb.Worksheets.Item("Asistencias")

# --- Asistencias: insert two new date columns before the old "Faltas" col (L) ---
$asistencias.Activate() | Out-Null
$asistencias.Range("L1:M1").EntireColumn.Insert()

# New exam date headers (same style as the other date headers in row 1)
$asistencias.Range("L1").Value = 45418
$asistencias.Range("M1").Value = 45425
$asistencias.Range("L1:M1").NumberFormat = $asistencias.Range("K1").NumberFormat
$asistencias.Range("L1:M1").HorizontalAlignment = $asistencias.Range("K1").HorizontalAlignment

# Match the column widths used for the new date columns
$asistencias.Columns("L:M").ColumnWidth = 6.6

# --- Concentrado: drop the split window, leave a plain selection at B4 ---
$concentrado.Activate() | Out-Null
$excel.ActiveWindow.FreezePanes = $false
$excel.ActiveWindow.Split = $false
$concentrado.Range("B4").Select() | Out-Null

# --- Asistencias becomes the active tab/sheet with K2 selected ---
$asistencias.Activate() | Out-Null
$asistencias.Range("K2").Select() | Out-Null
